$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '60.308.24'
Set-TextValue 'E2' '  -0.93%  '
Set-TextValue 'D3' '2.344.97'
Set-TextValue 'E3' '  -3.46%  '
Set-TextValue 'E4' '  +0.15%  '
Set-TextValue 'D5' '538.86'
Set-TextValue 'E5' '  -0.63%  '
Set-TextValue 'D6' '135.36'
Set-TextValue 'E6' '  -5.98%  '
Set-TextValue 'E7' '  +0.15%  '
Set-TextValue 'E8' '  -11.40%  '
Set-TextValue 'D9' '2.343.41'
Set-TextValue 'E9' '  -3.50%  '
Set-TextValue 'D10' '0.103'
Set-TextValue 'E10' '  -1.47%  '
Set-TextValue 'E11' '  +0.22%  '
Set-TextValue 'D12' '5.21'
Set-TextValue 'E12' '  -3.03%  '
Set-TextValue 'D13' '0.338'
Set-TextValue 'E13' '  -2.62%  '
Set-TextValue 'D14' '24.19'
Set-TextValue 'E14' '  -5.79%  '
Set-TextValue 'D15' '2.770.00'
Set-TextValue 'E15' '  -3.29%  '
Set-TextValue 'D16' '60.589.51'
Set-TextValue 'E16' '  -0.24%  '
Set-TextValue 'D17' '0.0000159'
Set-TextValue 'E17' '  -2.42%  '
Set-TextValue 'D18' '2.348.06'
Set-TextValue 'E18' '  -3.60%  '
Set-TextValue 'D19' '10.48'
Set-TextValue 'E19' '  -4.08%  '
Set-TextValue 'D20' '312.04'
Set-TextValue 'E20' '  -0.95%  '
Set-TextValue 'D21' '4.02'
Set-TextValue 'E21' '  -2.35%  '
Set-TextValue 'D22' '6.52'
Set-TextValue 'E22' '  -5.38%  '
Set-TextValue 'D23' '5.81'
Set-TextValue 'E23' '  -1.70%  '
Set-TextValue 'D24' '0.997'
Set-TextValue 'E24' '  -0.38%  '
Set-TextValue 'D25' '1.86'
Set-TextValue 'E25' '  +0.98%  '
Set-TextValue 'D26' '62.91'
Set-TextValue 'E26' '  -0.67%  '
Set-TextValue 'D27' '8.39'
Set-TextValue 'E27' '  +9.77%  '
Set-TextValue 'E28' '  +0.47%  '
Set-TextValue 'D29' '2.469.71'
Set-TextValue 'E29' '  -3.84%  '
Set-TextValue 'B30' 'InternetComputer(DFINITY)'
Set-TextValue 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D30' '7.83'
Set-TextValue 'E30' '  -3.65%  '
Set-TextValue 'B31' 'PEPE'
Set-TextValue 'C31' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D31' '0.0₃0876'
Set-TextValue 'E31' '  -7.73%  '
Set-TextValue 'B32' 'Bittensor'
Set-TextValue 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D32' '492.63'
Set-TextValue 'E32' '  -4.91%  '
Set-TextValue 'B33' 'Fetch.AI'
Set-TextValue 'C33' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D33' '1.35'
Set-TextValue 'E33' '  -6.67%  '
Set-TextValue 'E34' '  -1.66%  '
Set-TextValue 'E35' '  -6.29%  '
Set-TextValue 'D36' '1.50'
Set-TextValue 'E36' '  -3.69%  '
Set-TextValue 'D37' '1.00'
Set-TextValue 'E37' '  +0.12%  '
Set-TextValue 'B38' 'PolygonEcosystemToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D38' '0.370'
Set-TextValue 'E38' '  -0.64%  '
Set-TextValue 'B39' 'NEARProtocol'
Set-TextValue 'C39' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D39' '4.48'
Set-TextValue 'E39' '  -6.06%  '
Set-TextValue 'D40' '18.17'
Set-TextValue 'E40' '  -0.40%  '
Set-TextValue 'D41' '5.14'
Set-TextValue 'E41' '  -8.46%  '
Set-TextValue 'D42' '1.75'
Set-TextValue 'E42' '  +0.33%  '
Set-TextValue 'E43' '  -0.05%  '
Set-TextValue 'D44' '138.14'
Set-TextValue 'E44' '  -2.89%  '
Set-TextValue 'D45' '39.85'
Set-TextValue 'E45' '  -1.09%  '
Set-TextValue 'D46' '139.83'
Set-TextValue 'E46' '  +0.33%  '
Set-TextValue 'B47' 'Filecoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D47' '3.49'
Set-TextValue 'E47' '  -2.21%  '
Set-TextValue 'B48' 'dogwifhat'
Set-TextValue 'C48' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D48' '2.06'
Set-TextValue 'E48' '  -8.14%  '
Set-TextValue 'D49' '0.0504'
Set-TextValue 'E49' '  -4.44%  '
Set-TextValue 'D50' '19.11'
Set-TextValue 'E50' '  -9.40%  '
Set-TextValue 'E51' '  -2.75%  '
